$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing segment-name
# column (A) to B, PercActivations (B) to C, and PercSegmentAreas (C) to D.
$ws.Range("A1").EntireColumn.Insert()

# The segment-name column (now B) is plain data, so it should carry no special
# formatting (the header-style formatting that came along with the column
# insert needs to be stripped from the data rows).
$ws.Range("B2:B20").ClearFormats()

# New header for the index column, using the same header formatting as the
# other header cells (bold, centered, bordered).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# New column A: zero-based segment index, formatted like the header cells
# (this mirrors the style the segment-name column used to carry).
$ws.Range("C1").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)  # xlPasteFormats
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
}

$excel.CutCopyMode = $false
